# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# per the latest scrape. Mirrors the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.038.97'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.46%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.820.56'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.59%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -1.33%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.15'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.61%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.08%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4223'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.10%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3677'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07213'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.90%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8402'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.34%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.79'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.95%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.824.74'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.38%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.651'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.30%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07057'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.85%  '

# Row 15
$ws.Range("E15").Value = '  -3.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.86'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.74%  '

# Row 17
$ws.Range("E17").Value = '  -1.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008779'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.49%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.0000'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.12%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.92'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.99%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.105.76'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.23%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.129'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.79%  '

# Row 23
$ws.Range("E23").Value = '  -2.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.045.67'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.54%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.976'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.86%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.78'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.243'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.24'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.07%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.281'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.99%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.18'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.59%  '

# Row 31
$ws.Range("E31").Value = '  -2.25%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.176'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.59%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7389'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.14%  '

# Row 34
$ws.Range("E34").Value = '  -0.35%  '

# Row 36
$ws.Range("E36").Value = '  -1.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.091'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.93%  '

# Row 38
$ws.Range("E38").Value = '  -1.48%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05247'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.29%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.350'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.13%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.870'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.75%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1690'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.31%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5035'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.583'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.52'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.34%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.17'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.21%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4703'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.17%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.18%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06340'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.18%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.890'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.77%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.647'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.86%  '
